$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '62.826.58'
$ws.Range('E2').Value = '  -0.25%  '

# Row 3
$ws.Range('D3').Value = '2.466.74'
$ws.Range('E3').Value = '  +0.58%  '

# Row 4
$ws.Range('E4').Value = '  +0.05%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '574.55'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.55%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.71'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.06%  '

# Row 7
$ws.Range('E7').Value = '  -0.04%  '

# Row 8
$ws.Range('E8').Value = '  -1.11%  '

# Row 9
$ws.Range('D9').Value = '2.466.62'
$ws.Range('E9').Value = '  +0.59%  '

# Row 10
$ws.Range('E10').Value = '  -0.32%  '

# Row 11
$ws.Range('E11').Value = '  -0.34%  '

# Row 12
$ws.Range('E12').Value = '  -0.38%  '

# Row 13
$ws.Range('E13').Value = '  +0.57%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '29.31'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +3.18%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.0000178'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -1.73%  '

# Row 16
$ws.Range('D16').Value = '2.915.10'
$ws.Range('E16').Value = '  +0.71%  '

# Row 17
$ws.Range('D17').Value = '62.748.60'
$ws.Range('E17').Value = '  -0.03%  '

# Row 18
$ws.Range('D18').Value = '2.482.71'
$ws.Range('E18').Value = '  +1.21%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.96'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.19%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.00'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.66%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '326.87'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -1.88%  '

# Row 22
$ws.Range('E22').Value = '  -0.05%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.21'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +7.02%  '

# Row 24
$ws.Range('E24').Value = '  +0.00%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.07'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +17.89%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '65.60'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.25%  '

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '641.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -1.75%  '

# Row 28
$ws.Range('D28').Value = '0.0₃0983'
$ws.Range('E28').Value = '  -1.87%  '

# Row 29
$ws.Range('D29').Value = '2.591.45'
$ws.Range('E29').Value = '  +0.71%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.996'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -15.37%  '

# Row 31
$ws.Range('E31').Value = '  -0.74%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.94'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.44%  '

# Row 33
$ws.Range('E33').Value = '  -2.11%  '

# Row 34
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.135'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.25%  '

# Row 35
$ws.Range('E35').Value = '  -0.04%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.55'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.33%  '

# Row 37
$ws.Range('E37').Value = '  -0.56%  '

# Row 38
$ws.Range('B38').Value = 'Monero'
$ws.Range('C38').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '152.19'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -0.56%  '

# Row 39
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.369'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -1.57%  '

# Row 40
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.80'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.53%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '18.68'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.77%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.37'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -3.35%  '

# Row 43
$ws.Range('E43').Value = '  -2.08%  '

# Row 45
$ws.Range('E45').Value = '  -29.74%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '152.79'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +4.81%  '

# Row 47
$ws.Range('E47').Value = '  +1.68%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '3.59'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.62%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '20.48'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.10%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.608'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +0.33%  '

# Row 51
$ws.Range('E51').Value = '  -1.35%  '
